$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 3
$ws.Range("E3").Value = 24
$ws.Range("F3").Value = 12
$ws.Range("H3").Value = 12

# Row 15
$ws.Range("E15").Value = 151
$ws.Range("F15").Value = 79
$ws.Range("H15").Value = 79

# Row 17
$ws.Range("E17").Value = 103
$ws.Range("F17").Value = 47
$ws.Range("H17").Value = 47

# Row 20
$ws.Range("E20").Value = 6

# Row 27
$ws.Range("E27").Value = 8

# Row 39
$ws.Range("E39").Value = 24

# Row 41
$ws.Range("E41").Value = 35

# Row 46
$ws.Range("E46").Value = 26
$ws.Range("F46").Value = 9
$ws.Range("H46").Value = 9

# Row 47
$ws.Range("E47").Value = 56
$ws.Range("F47").Value = 33
$ws.Range("H47").Value = 33

# Row 49
$ws.Range("E49").Value = 61

# Row 63
$ws.Range("E63").Value = 31

# Row 66
$ws.Range("E66").Value = 32
$ws.Range("F66").Value = 19
$ws.Range("H66").Value = 19

# Row 80
$ws.Range("E80").Value = 22
